$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 285; this shifts the existing rows
# 285-337 down to 286-338 (matches the dimension growing to A1:R338).
$ws.Rows.Item(285).Insert()

# Populate the newly inserted row 285 with the new record's data.
$ws.Cells.Item(285, 1).Value = 3
$ws.Cells.Item(285, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(285, 3).Value = "Coquimbo"
$ws.Cells.Item(285, 4).Value = 44694
$ws.Cells.Item(285, 5).Value = 5
$ws.Cells.Item(285, 6).Value = 100112012
$ws.Cells.Item(285, 7).Value = "Espinaca"
$ws.Cells.Item(285, 8).Value = "Sin especificar"
$ws.Cells.Item(285, 9).Value = "Primera"
$ws.Cells.Item(285, 10).Value = 175
$ws.Cells.Item(285, 11).Value = 3500
$ws.Cells.Item(285, 12).Value = 4000
$ws.Cells.Item(285, 13).Value = 3743
$ws.Cells.Item(285, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(285, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(285, 16).Value = 1248
$ws.Cells.Item(285, 17).Value = 3
$ws.Cells.Item(285, 18).Value = "Hortaliza"
